$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = 636

$ws.Range("A27").Value = "Örebro University"
$ws.Range("B27").Value = 981
$ws.Range("C27").Value = 917
$ws.Range("A28").Value = "Örebro University Hospital"
$ws.Range("B28").Value = 182
$ws.Range("C28").Value = 0
$ws.Range("A29").Value = "Oslo University Hospital"
$ws.Range("B29").Value = 740
$ws.Range("C29").Value = 1057
$ws.Range("A30").Value = "Oulu University Hospital"
$ws.Range("B30").Value = 1213
$ws.Range("C30").ClearContents()
$ws.Range("A31").Value = "Sahlgrenska University Hospital"
$ws.Range("B31").Value = 867
$ws.Range("C31").Value = 1201
$ws.Range("A32").Value = "Skane University Hospital"
$ws.Range("B32").Value = 523
$ws.Range("C32").Value = 1119
$ws.Range("A33").Value = "St. Olav’s University Hospital"
$ws.Range("B33").Value = 1331
$ws.Range("C33").ClearContents()
$ws.Range("A34").Value = "Steno Diabetes Center Copenhagen"
$ws.Range("B34").Value = 560
$ws.Range("C34").Value = 408.0000000000001
$ws.Range("A35").Value = "Stockholm South General Hospital"
$ws.Range("B35").Value = 216
$ws.Range("C35").Value = 1052
$ws.Range("A36").Value = "Tampere University Hospital"
$ws.Range("B36").Value = 400.5000000000001
$ws.Range("C36").Value = 734.0000000000001
$ws.Range("A37").Value = "The National University Hospital of Iceland"
$ws.Range("B37").Value = 992
$ws.Range("C37").Value = 765.9999999999999
$ws.Range("A38").Value = "Turku University Hospital"
$ws.Range("B38").Value = 913
$ws.Range("C38").ClearContents()
$ws.Range("A39").Value = "UiT The Arctic University of Norway"
$ws.Range("B39").Value = 655.4999999999999
$ws.Range("C39").ClearContents()
$ws.Range("A40").Value = "Umeå University"
$ws.Range("B40").Value = 909.9999999999999
$ws.Range("C40").ClearContents()
$ws.Range("A41").Value = "University Hospital of North Norway"
$ws.Range("B41").Value = 973
$ws.Range("C41").ClearContents()
$ws.Range("A42").Value = "University Hospital of Umeå"
$ws.Range("B42").Value = 427
$ws.Range("C42").Value = 609.9999999999999
$ws.Range("A43").Value = "University of Bergen"
$ws.Range("B43").Value = 689
$ws.Range("C43").ClearContents()
$ws.Range("A44").Value = "University of Copenhagen"
$ws.Range("B44").Value = 846
$ws.Range("C44").Value = 1220
$ws.Range("A45").Value = "University of Eastern Finland"
$ws.Range("B45").Value = 1278
$ws.Range("C45").Value = 1292
$ws.Range("A46").Value = "University of Helsinki"
$ws.Range("B46").Value = 723
$ws.Range("C46").Value = 1163
$ws.Range("A47").Value = "University of Iceland"
$ws.Range("B47").Value = 1306
$ws.Range("C47").Value = 828.9999999999998
$ws.Range("A48").Value = "University of Oslo"
$ws.Range("B48").Value = 670
$ws.Range("C48").Value = 566
$ws.Range("A49").Value = "University of Oulu"
$ws.Range("B49").Value = 1191
$ws.Range("C49").Value = 1832
$ws.Range("A50").Value = "University of Southern Denmark"
$ws.Range("B50").Value = 577.4999999999999
$ws.Range("C50").Value = 510.9999999999999
$ws.Range("A51").Value = "University of Tampere"
$ws.Range("B51").Value = 276
$ws.Range("C51").Value = 701.9999999999999
$ws.Range("A52").Value = "University of Turku"
$ws.Range("B52").Value = 759
$ws.Range("C52").Value = 1194.5
$ws.Range("A53").Value = "Uppsala Academic Hospital"
$ws.Range("B53").Value = 948.9999999999999
$ws.Range("C53").Value = 594.9999999999999
$ws.Range("A54").Value = "Uppsala University"
$ws.Range("B54").Value = 971
$ws.Range("C54").ClearContents()
$ws.Range("A55").Value = "Zealand University Hospital"
$ws.Range("B55").Value = 593.5
$ws.Range("C55").Value = 541.5

$ws.Range("B57").Value = 720
$ws.Range("C60").Value = 1217
$ws.Range("B61").Value = 697.5000000000001
$ws.Range("C61").Value = 1123.5
